$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 338

# Update row 3 values
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 118

# Remove rows 4 and 5 (shift cells up), shrinking the used range to A1:B3
$ws.Range("A4:B5").Delete()
